# Regenerate the experiment task-order sheets (new randomized CSV filenames)
# and re-label/re-order the task-order tabs accordingly.
$wb = $excel.ActiveWorkbook

# --- Sheet 1 (was GNG_TO) becomes vSAT_TO: same 4-row shape, new filenames ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "vSAT_stims-1651589014423244.csv"
$ws1.Range("B3").Value = "SAT_stims-1651589014391993.csv"
$ws1.Range("B4").Value = "SAT_stims-16515890144076183.csv"
$ws1.Range("B5").Value = "vSAT_stims-16515890144388692.csv"
$ws1.Name = "vSAT_TO-16515890144544935"

# --- Sheet 2 (was NB_TO, 10 rows) becomes RS_TO: shrinks to the 3-row
#     "eyes closed"/"eyes open" resting-state shape ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows("4:10").Delete()
$ws2.Range("B2").Value = "eyes closed"
$ws2.Range("B3").Value = "eyes open"
$ws2.Name = "RS_TO-16515890144544935"

# --- Sheet 3 (was RS_TO, 3 rows) becomes GNG_TO: grows to the 5-row
#     go/no-go filename shape ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows("4:5").Insert()
$ws3.Range("A2").Copy()
$ws3.Range("A4:A5").PasteSpecial(-4122)
$ws3.Range("A4").Value = 2
$ws3.Range("A5").Value = 3
$ws3.Range("B2").Value = "go_stims-16515890144544935.csv"
$ws3.Range("B3").Value = "GNG_stims-1651589014470118.csv"
$ws3.Range("B4").Value = "go_stims-1651589014470118.csv"
$ws3.Range("B5").Value = "GNG_stims-1651589014485743.csv"
$ws3.Name = "GNG_TO-1651589014485743"

# --- Sheet 4 (TOL_TO, 7 rows) stays TOL_TO: same shape, new filenames ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16515890145013683.csv"
$ws4.Range("B3").Value = "ZM_stims-1651589014485743.csv"
$ws4.Range("B4").Value = "MM_stims-1651589014516994.csv"
$ws4.Range("B5").Value = "ZM_stims-16515890145013683.csv"
$ws4.Range("B6").Value = "MM_stims-16515890145326183.csv"
$ws4.Range("B7").Value = "ZM_stims-1651589014516994.csv"
$ws4.Name = "TOL_TO-16515890145326183"

# --- Sheet 5 (was vSAT_TO, 5 rows) becomes NB_TO: grows to the 10-row
#     n-back filename shape ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Rows("6:10").Insert()
$ws5.Range("A2").Copy()
$ws5.Range("A6:A10").PasteSpecial(-4122)
$ws5.Range("A6").Value = 4
$ws5.Range("A7").Value = 5
$ws5.Range("A8").Value = 6
$ws5.Range("A9").Value = 7
$ws5.Range("A10").Value = 8
$ws5.Range("B2").Value = "TB-16515890160841186.csv"
$ws5.Range("B3").Value = "OB-16515890160060296.csv"
$ws5.Range("B4").Value = "TB-16515890166371741.csv"
$ws5.Range("B5").Value = "ZB-match_3-16515890147625074.csv"
$ws5.Range("B6").Value = "OB-16515890149330254.csv"
$ws5.Range("B7").Value = "ZB-match_3-16515890146353676.csv"
$ws5.Range("B8").Value = "OB-1651589015378874.csv"
$ws5.Range("B9").Value = "ZB-match_6-16515890147868676.csv"
$ws5.Range("B10").Value = "TB-16515890167152994.csv"
$ws5.Name = "NB_TO-16515890167309248"
